# Add payment 71277620 (Cash) 2025-08-18T17:04:15
# 1. Fix existing row 29, column A (phone) to be stored as a number instead of text.
# 2. Append a new row 30 with the new payment record (phone stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: correct A29 from text "71277620" to the numeric value 71277620 ---
$ws.Cells.Item(29, 1).Value = 71277620

# --- Row 30: new payment record ---
# A30 ("phone") must stay text, even though it looks numeric, so force a text
# number format before assigning the value, then drop the format again so the
# cell keeps the default style (matching the style-less cells around it).
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "71277620"
$ws.Cells.Item(30, 1).ClearFormats()

# B30 ("amount") is left blank, same as the source record.
$ws.Cells.Item(30, 3).Value = "Cash"
$ws.Cells.Item(30, 4).Value = "2025-08-18T17:04:15"
$ws.Cells.Item(30, 5).Value = 760
# F30 ("discount_applied") is left blank, same as the source record.
$ws.Cells.Item(30, 7).Value = 684
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 76
